# Updated symbol list on Fri Dec 23 15:25:05 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.34"
$ws.Range("D3").Value = "'21.94"
$ws.Range("D4").Value = "'5.375"
$ws.Range("D5").Value = "'0.05855"
$ws.Range("D6").Value = "'3.388"
$ws.Range("D7").Value = "'6.364"
$ws.Range("D8").Value = "'0.8164"
$ws.Range("D9").Value = "'1.009"
$ws.Range("D10").Value = "'0.01115"
$ws.Range("D12").Value = "'0.03706"
$ws.Range("D13").Value = "'0.07425"
$ws.Range("D14").Value = "'0.03036"
$ws.Range("D15").Value = "'4.189"
$ws.Range("D16").Value = "'0.09395"
$ws.Range("D17").Value = "'0.001594"
$ws.Range("D18").Value = "'0.04828"
$ws.Range("D19").Value = "'0.006050"
$ws.Range("D20").Value = "'0.004091"
$ws.Range("D21").Value = "'0.0009897"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D23").Value = "'3.700"
$ws.Range("D24").Value = "'2.224"
$ws.Range("D25").Value = "'0.3238"
$ws.Range("D26").Value = "'0.1295"
$ws.Range("D27").Value = "'0.0002492"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("D40").Value = "'0.03860"
$ws.Range("D41").Value = "'0.006435"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.002598"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").Value = "'0.006231"
$ws.Range("D45").Value = "'0.00005618"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.6995"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.01009"
